# GANTT.xlsx edit: "commit cahier des charges + cahier de recette"
#
# Swaps / updates a couple of task rows, adds a new task row (overwriting the
# stray "Activity 10" label that used to sit in row 14), adjusts plan/actual
# periods accordingly, and moves the current selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: was "Créer une page pour afficher ... l'utilisateur" -------------
$ws.Range("B6").Value = "Faire le design du site"
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 1

# --- Row 7: was "Faire la design du site" ------------------------------------
$ws.Range("B7").Value = "Créer une page pour afficher la position et la vitesse de l'utilisateur"
$ws.Range("C7").Value = 11
$ws.Range("D7").Value = 12
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0

# --- Row 5: plan/actual periods updated --------------------------------------
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 1

# --- Row 13: "Faire la bdd" actual period filled in --------------------------
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = 1

# --- Row 14: used to show the orphan "Activity 10" label, now becomes a new --
# --- task row: "Créer la page pour les administrateurs" ----------------------
$ws.Range("B14").Value = "Créer la page pour les administrateurs"
$ws.Range("C14").Value = 11
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 11
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 1
$ws.Rows.Item(14).RowHeight = 96.75

# --- view: scrolled down a bit, selection on E13 -----------------------------
$ws.Application.Goto($ws.Range("A10"), $false)
$ws.Range("E13").Select()
